$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues
$scratch = $ws.Range("ZZ1")

function Set-TextValue {
    param([string]$addr, [string]$value)
    $formulaLiteral = $value.Replace('"', '""')
    $scratch.Formula = '="' + $formulaLiteral + '"'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial($xlPasteValues)
}

Set-TextValue "D2" '43.498.33'
Set-TextValue "E2" '  +2.86%  '
Set-TextValue "D3" '2.313.66'
Set-TextValue "E3" '  +1.88%  '
Set-TextValue "E4" '  +0.01%  '
Set-TextValue "D5" '310.68'
Set-TextValue "E5" '  +0.67%  '
Set-TextValue "D6" '105.03'
Set-TextValue "E6" '  +7.72%  '
Set-TextValue "E7" '  +1.56%  '
Set-TextValue "E8" '  +0.06%  '
Set-TextValue "D10" '36.92'
Set-TextValue "E10" '  +5.56%  '
Set-TextValue "D11" '52.90'
Set-TextValue "D12" '0.0815'
Set-TextValue "E12" '  +0.56%  '
Set-TextValue "E13" '  -0.88%  '
Set-TextValue "E14" '  +2.72%  '
Set-TextValue "D15" '2.672.92'
Set-TextValue "E15" '  +1.88%  '
Set-TextValue "D16" '15.15'
Set-TextValue "E16" '  +3.97%  '
Set-TextValue "D17" '2.305.56'
Set-TextValue "E17" '  +1.93%  '
Set-TextValue "E18" '  +3.43%  '
Set-TextValue "D19" '43.411.05'
Set-TextValue "E19" '  +2.95%  '
Set-TextValue "D20" '12.24'
Set-TextValue "E20" '  -0.20%  '
Set-TextValue "D21" '0.0₃0931'
Set-TextValue "E21" '  +2.92%  '
Set-TextValue "D22" '6.20'
Set-TextValue "E22" '  +4.00%  '
Set-TextValue "D23" '68.35'
Set-TextValue "E23" '  +1.11%  '
Set-TextValue "D24" '242.94'
Set-TextValue "E24" '  +2.71%  '
Set-TextValue "E25" '  +3.44%  '
Set-TextValue "E26" '  +0.56%  '
Set-TextValue "E27" '  +0.02%  '
Set-TextValue "D28" '24.86'
Set-TextValue "E28" '  +5.45%  '
Set-TextValue "D30" '37.10'
Set-TextValue "E30" '  -0.06%  '
Set-TextValue "D31" '9.67'
Set-TextValue "E31" '  +1.18%  '
Set-TextValue "D32" '166.55'
Set-TextValue "E32" '  +2.05%  '
Set-TextValue "D33" '5.31'
Set-TextValue "E33" '  +1.30%  '
Set-TextValue "E34" '  +0.00%  '
Set-TextValue "D35" '18.42'
Set-TextValue "E35" '  +4.89%  '
Set-TextValue "E36" '  +6.96%  '
Set-TextValue "D37" '0.0747'
Set-TextValue "E37" '  +2.05%  '
Set-TextValue "E38" '  -0.97%  '
Set-TextValue "D39" '4.60'
Set-TextValue "E39" '  +10.71%  '
Set-TextValue "E40" '  +3.79%  '
Set-TextValue "D41" '0.107'
Set-TextValue "E41" '  +2.79%  '
Set-TextValue "E42" '  +0.70%  '
Set-TextValue "E43" '  +20.02%  '
Set-TextValue "E44" '  +3.96%  '
Set-TextValue "D45" '1.997.39'
Set-TextValue "E45" '  +2.62%  '
Set-TextValue "D46" '19.13'
Set-TextValue "E46" '  +1.31%  '
Set-TextValue "D47" '3.14'
Set-TextValue "E47" '  +6.83%  '
Set-TextValue "D48" '10.04'
Set-TextValue "E48" '  +2.70%  '
Set-TextValue "D49" '57.04'
Set-TextValue "E49" '  +4.91%  '
Set-TextValue "D50" '2.95'
Set-TextValue "E50" '  +1.49%  '

$scratch.Clear()
$excel.CutCopyMode = $false

